$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit new longer content (best-fit re-computed by Excel
# after the longer strings below are entered)
$ws.Columns.Item(2).ColumnWidth = 41

# Row 17: Work On Design
$ws.Range("B17").Value = "Work On Design"

# Row 20: Between Phases there will be a short cutscene
$ws.Range("B20").Value = "Between Phases there will be a short cutscene"

# Row 22: Define Phase 2 (moved from old B18)
$ws.Range("B22").Value = "Define Phase 2"

# Row 23: Create second Scene
$ws.Range("B23").Value = "Create second Scene"

# Row 24: Jumping Puzzle
$ws.Range("B24").Value = "Jumping Puzzle"

# Row 25: Re-design-Boss
$ws.Range("B25").Value = "Re-design-Boss"

# Row 26: Add Orbs (same as first fight)
$ws.Range("B26").Value = "Add Orbs (same as first fight)"

# Row 27: add Floor wipe mechanic
$ws.Range("B27").Value = "add Floor wipe mechanic"

# Row 28: Work on the winning variables () / still using the one from Level One
$ws.Range("B28").Value = "Work on the winning variables ()"
$ws.Range("C28").Value = "still using the one from Level One"

# Remove old row 18 content (the value was moved to B22 above)
$ws.Range("B18").ClearContents()

# Update selection to match target workbook view
$ws.Range("C31").Select()
